# Insert a new data row at row 138 (pushing the existing rows 138:150 down
# to 139:151) and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 138:150 down to 139:151, leaving a blank row 138 behind.
$ws.Rows.Item(138).Insert()

# Fill in the new row 138 with the new "Choclo" record.
$ws.Cells.Item(138, 1).Value = 4
$ws.Cells.Item(138, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(138, 3).Value = "Los Lagos"
$ws.Cells.Item(138, 4).Value = 44491
$ws.Cells.Item(138, 5).Value = 10
$ws.Cells.Item(138, 6).Value = 100112024
$ws.Cells.Item(138, 7).Value = "Choclo"
$ws.Cells.Item(138, 8).Value = "Dulce o Americano"
$ws.Cells.Item(138, 9).Value = "Primera"
$ws.Cells.Item(138, 10).Value = 100
$ws.Cells.Item(138, 11).Value = 40000
$ws.Cells.Item(138, 12).Value = 40000
$ws.Cells.Item(138, 13).Value = 40000
$ws.Cells.Item(138, 14).Value = "`$/malla 70 unidades"
$ws.Cells.Item(138, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(138, 16).Value = 571
$ws.Cells.Item(138, 17).Value = 70
$ws.Cells.Item(138, 18).Value = "Hortaliza"
